# Apply cryptos list update (price/volume refresh + two row-label swaps)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text cell updates (coin names, links, and Volume/% cells, plus multi-dot
# "Price" strings that Excel cannot misinterpret as numbers).
$plainUpdates = @{
    'D2' = '27.699.63'
    'E2' = '  -0.11%  '
    'D3' = '1.845.79'
    'E3' = '  -0.88%  '
    'E4' = '  -2.43%  '
    'E5' = '  -1.32%  '
    'E6' = '  -2.03%  '
    'E7' = '  -2.68%  '
    'E8' = '  -1.77%  '
    'E9' = '  -1.60%  '
    'E10' = '  -1.27%  '
    'E11' = '  -0.47%  '
    'D12' = '1.851.47'
    'E12' = '  -1.41%  '
    'E13' = '  -0.74%  '
    'E14' = '  -1.85%  '
    'E15' = '  -1.14%  '
    'E16' = '  +5.05%  '
    'E18' = '  -1.67%  '
    'E20' = '  -0.12%  '
    'D21' = '27.683.55'
    'E21' = '  -0.24%  '
    'E22' = '  -1.96%  '
    'E23' = '  -1.85%  '
    'D24' = '2.074.05'
    'E24' = '  -0.97%  '
    'E26' = '  -1.74%  '
    'E27' = '  -1.40%  '
    'E28' = '  +9.23%  '
    'E29' = '  -0.13%  '
    'E30' = '  -0.07%  '
    'E31' = '  -1.02%  '
    'E32' = '  -0.19%  '
    'E33' = '  -0.51%  '
    'E34' = '  -1.13%  '
    'E35' = '  -4.50%  '
    'E37' = '  -1.06%  '
    'E38' = '  -0.30%  '
    'E39' = '  -0.90%  '
    'E40' = '  +5.55%  '
    'E41' = '  +1.13%  '
    'B42' = 'TheSandbox'
    'C42' = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
    'E42' = '  -1.88%  '
    'B43' = 'Algorand'
    'C43' = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
    'E43' = '  -0.50%  '
    'E44' = '  +0.85%  '
    'B45' = 'Quant'
    'C45' = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
    'E45' = '  -1.89%  '
    'B46' = 'EnergySwap'
    'C46' = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
    'E46' = '  -1.22%  '
    'E47' = '  +0.34%  '
    'E48' = '  -3.00%  '
    'E49' = '  -1.99%  '
    'E50' = '  -1.18%  '
    'E51' = '  -3.73%  '
}
foreach ($addr in $plainUpdates.Keys) {
    $ws.Range($addr).Value = $plainUpdates[$addr]
}

# "Price" cells whose new text looks like a plain number (e.g. "319.11").
# These must stay plain text (as in the original workbook), so force the cell
# to Text format before assigning, then restore the default "Normal" style so
# no extraneous formatting is left behind on the cell.
$numericLookingUpdates = @{
    'D4' = '1.011'
    'D5' = '319.11'
    'D6' = '1.010'
    'D7' = '0.4308'
    'D8' = '0.3740'
    'D9' = '0.07341'
    'D10' = '0.8766'
    'D11' = '21.59'
    'D13' = '6.723'
    'D14' = '5.441'
    'D15' = '0.07117'
    'D16' = '88.70'
    'D17' = '1.014'
    'D18' = '0.000008975'
    'D20' = '15.46'
    'D22' = '5.210'
    'D23' = '11.10'
    'D26' = '155.55'
    'D27' = '18.62'
    'D28' = '2.167'
    'D29' = '5.367'
    'D30' = '118.89'
    'D31' = '0.08937'
    'D32' = '1.231'
    'D33' = '0.7754'
    'D34' = '4.549'
    'D35' = '2.884'
    'D36' = '1.012'
    'D37' = '1.133'
    'D38' = '0.05334'
    'D39' = '0.01972'
    'D40' = '7.287'
    'D41' = '2.923'
    'D42' = '0.5114'
    'D43' = '0.1683'
    'D44' = '8.795'
    'D45' = '109.10'
    'D46' = '10.65'
    'D47' = '0.4745'
    'D48' = '0.06471'
    'D49' = '1.012'
    'D50' = '1.691'
    'D51' = '1.846'
}
foreach ($addr in $numericLookingUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $numericLookingUpdates[$addr]
    $cell.Style = "Normal"
}

